# Sravana poornima 2019 - update donors list amounts to formatted INR currency text
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$amounts = @{
    2 = "₹ 3,50,000"
    3 = "₹ 3,30,000"
    4 = "₹ 3,00,000"
    5 = "₹ 1,50,000"
    6 = "₹ 1,00,116"
    7 = "₹ 1,00,000"
    8 = "₹ 1,00,000"
    9 = "₹ 1,00,000"
    10 = "₹ 1,00,000"
    11 = "₹ 1,00,000"
    12 = "₹ 1,00,000"
    13 = "₹ 1,00,000"
    14 = "₹ 1,00,000"
    15 = "₹ 60,000"
    16 = "₹ 50,116"
    17 = "₹ 50,116"
    18 = "₹ 50,011"
    19 = "₹ 50,000"
    20 = "₹ 50,000"
    21 = "₹ 50,000"
    22 = "₹ 50,000"
    23 = "₹ 50,000"
    24 = "₹ 50,000"
    25 = "₹ 50,000"
    26 = "₹ 50,000"
    27 = "₹ 50,000"
    28 = "₹ 50,000"
    29 = "₹ 50,000"
    30 = "₹ 32,500"
    31 = "₹ 25,000"
    32 = "₹ 25,000"
    33 = "₹ 25,000"
    34 = "₹ 25,000"
    35 = "₹ 25,000"
    36 = "₹ 25,000"
    37 = "₹ 25,000"
    38 = "₹ 20,232"
    39 = "₹ 20,000"
    40 = "₹ 20,000"
    41 = "₹ 20,000"
    42 = "₹ 20,000"
    43 = "₹ 20,000"
    44 = "₹ 20,000"
    45 = "₹ 20,000"
    46 = "₹ 20,000"
    47 = "₹ 20,000"
    48 = "₹ 16,116"
    49 = "₹ 15,000"
    50 = "₹ 11,493"
    51 = "₹ 10,116"
    52 = "₹ 10,116"
    53 = "₹ 10,116"
    54 = "₹ 10,116"
    55 = "₹ 10,116"
    56 = "₹ 10,116"
    57 = "₹ 10,116"
    58 = "₹ 10,116"
    59 = "₹ 10,116"
    60 = "₹ 10,116"
    61 = "₹ 10,116"
    62 = "₹ 10,116"
    63 = "₹ 10,116"
    64 = "₹ 10,116"
    65 = "₹ 10,016"
    66 = "₹ 10,001"
    67 = "₹ 10,001"
    68 = "₹ 10,001"
    69 = "₹ 10,001"
    70 = "₹ 10,000"
    71 = "₹ 10,000"
    72 = "₹ 10,000"
    73 = "₹ 10,000"
    74 = "₹ 10,000"
    75 = "₹ 10,000"
    76 = "₹ 10,000"
    77 = "₹ 10,000"
    78 = "₹ 10,000"
    79 = "₹ 10,000"
    80 = "₹ 10,000"
    81 = "₹ 10,000"
    82 = "₹ 10,000"
    83 = "₹ 10,000"
    84 = "₹ 10,000"
    85 = "₹ 10,000"
    86 = "₹ 10,000"
    87 = "₹ 10,000"
    88 = "₹ 10,000"
    89 = "₹ 10,000"
    90 = "₹ 10,000"
    91 = "₹ 10,000"
    92 = "₹ 10,000"
    93 = "₹ 10,000"
    94 = "₹ 10,000"
    95 = "₹ 10,000"
    96 = "₹ 10,000"
    97 = "₹ 10,000"
    98 = "₹ 10,000"
    99 = "₹ 10,000"
    100 = "₹ 10,000"
    101 = "₹ 10,000"
    102 = "₹ 10,000"
    103 = "₹ 10,000"
    104 = "₹ 10,000"
    105 = "₹ 10,000"
    106 = "₹ 10,000"
    107 = "₹ 10,000"
    108 = "₹ 10,000"
    109 = "₹ 10,000"
    110 = "₹ 10,000"
    111 = "₹ 10,000"
    112 = "₹ 10,000"
    113 = "₹ 10,000"
    114 = "₹ 10,000"
    115 = "₹ 10,000"
    116 = "₹ 7,722"
    117 = "₹ 7,000"
    118 = "₹ 7,000"
    119 = "₹ 6,000"
    120 = "₹ 5,116"
    121 = "₹ 5,116"
    122 = "₹ 5,116"
    123 = "₹ 5,116"
    124 = "₹ 5,116"
    125 = "₹ 5,116"
    126 = "₹ 5,116"
    127 = "₹ 5,116"
    128 = "₹ 5,116"
    129 = "₹ 5,100"
    130 = "₹ 5,004"
    131 = "₹ 5,004"
    132 = "₹ 5,001"
    133 = "₹ 5,001"
    134 = "₹ 5,001"
    135 = "₹ 5,001"
    136 = "₹ 5,000"
    137 = "₹ 5,000"
    138 = "₹ 5,000"
    139 = "₹ 5,000"
    140 = "₹ 5,000"
    141 = "₹ 5,000"
    142 = "₹ 5,000"
    143 = "₹ 5,000"
    144 = "₹ 5,000"
    145 = "₹ 5,000"
    146 = "₹ 5,000"
    147 = "₹ 5,000"
    148 = "₹ 5,000"
    149 = "₹ 5,000"
    150 = "₹ 5,000"
    151 = "₹ 5,000"
    152 = "₹ 5,000"
    153 = "₹ 5,000"
    154 = "₹ 5,000"
    155 = "₹ 5,000"
    156 = "₹ 5,000"
    157 = "₹ 5,000"
    158 = "₹ 5,000"
    159 = "₹ 5,000"
    160 = "₹ 5,000"
    161 = "₹ 5,000"
    162 = "₹ 5,000"
    163 = "₹ 5,000"
    164 = "₹ 5,000"
    165 = "₹ 5,000"
    166 = "₹ 5,000"
    167 = "₹ 5,000"
    168 = "₹ 5,000"
    169 = "₹ 5,000"
    170 = "₹ 3,001"
    171 = "₹ 3,000"
    172 = "₹ 3,000"
    173 = "₹ 3,000"
    174 = "₹ 3,000"
    175 = "₹ 2,500"
    176 = "₹ 2,016"
    177 = "₹ 2,000"
    178 = "₹ 2,000"
    179 = "₹ 2,000"
    180 = "₹ 2,000"
    181 = "₹ 2,000"
    182 = "₹ 2,000"
    183 = "₹ 1,120"
    184 = "₹ 1,116"
    185 = "₹ 1,116"
    186 = "₹ 1,116"
    187 = "₹ 1,116"
    188 = "₹ 1,116"
    189 = "₹ 1,016"
    190 = "₹ 1,011"
    191 = "₹ 1,001"
    192 = "₹ 1,001"
    193 = "₹ 1,001"
    194 = "₹ 1,001"
    195 = "₹ 1,000"
    196 = "₹ 1,000"
    197 = "₹ 1,000"
    198 = "₹ 520"
    199 = "₹ 501"
    200 = "₹ 120"
}

foreach ($row in $amounts.Keys) {
    $ws.Cells.Item($row, 6).Value = $amounts[$row]
}

$ws.Range("F1:F1048576").Select()
